$d = $word.ActiveDocument

$oldText1 = "Results of your t-tests in table form for both experiments from Step 6"
$newText1 = "Results of your t-tests for both experiments from Step 6"

$oldText2 = "Explain why you rejected or failed to reject your hypotheses and explain the processes that influenced both experiments."
$newText2 = "Explain why you rejected or failed to reject your null hypotheses based on your t-test results and explain the processes that influenced both experiments."

foreach ($p in $d.Paragraphs) {
    $t = $p.Range.Text
    if ($t -eq ($oldText1 + "`r")) {
        $p.Range.Text = $newText1
    }
    elseif ($t -eq ($oldText2 + "`r")) {
        $p.Range.Text = $newText2
    }
}
